$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the bordered/centered header style (used by column A, s=1) down to the new rows
$ws.Range("A10").Copy()
$ws.Range("A11:A17").PasteSpecial(-4122)

# Row 2
$ws.Cells.Item(2, 1).Value = 34
$ws.Cells.Item(2, 2).Value = 35
$ws.Cells.Item(2, 3).Value = 364.1653337305623
$ws.Cells.Item(2, 4).Value = 0.03055558
$ws.Cells.Item(2, 5).Value = 0.03055558
$ws.Cells.Item(2, 6).Value = 0.06111116
$ws.Cells.Item(2, 7).Value = 546.2480005958435
$ws.Cells.Item(2, 8).Value = 0.38
$ws.Cells.Item(2, 9).Value = 0.43
$ws.Cells.Item(2, 10).Value = 0.11
$ws.Cells.Item(2, 11).Value = 0.37

# Row 3
$ws.Cells.Item(3, 1).Value = 35
$ws.Cells.Item(3, 2).Value = 36
$ws.Cells.Item(3, 3).Value = 381.7906987705871
$ws.Cells.Item(3, 4).Value = 0.05833338
$ws.Cells.Item(3, 5).Value = 0.05833338
$ws.Cells.Item(3, 6).Value = 0.11666676
$ws.Cells.Item(3, 7).Value = 572.6860481558806
$ws.Cells.Item(3, 8).Value = 0.37
$ws.Cells.Item(3, 9).Value = 0.49
$ws.Cells.Item(3, 10).Value = 0.12
$ws.Cells.Item(3, 11).Value = 0.3

# Row 4
$ws.Cells.Item(4, 1).Value = 36
$ws.Cells.Item(4, 2).Value = 37
$ws.Cells.Item(4, 3).Value = 410.3627591381613
$ws.Cells.Item(4, 4).Value = 0.07777784
$ws.Cells.Item(4, 5).Value = 0.07777784
$ws.Cells.Item(4, 6).Value = 0.15555568
$ws.Cells.Item(4, 7).Value = 615.5441387072419
$ws.Cells.Item(4, 8).Value = 0.4
$ws.Cells.Item(4, 9).Value = 0.46
$ws.Cells.Item(4, 10).Value = 0.19
$ws.Cells.Item(4, 11).Value = 0.33

# Row 5
$ws.Cells.Item(5, 1).Value = 37
$ws.Cells.Item(5, 2).Value = 38
$ws.Cells.Item(5, 3).Value = 408.4460994869478
$ws.Cells.Item(5, 4).Value = 0.07222228
$ws.Cells.Item(5, 5).Value = 0.07222228
$ws.Cells.Item(5, 6).Value = 0.14444456
$ws.Cells.Item(5, 7).Value = 612.6691492304217
$ws.Cells.Item(5, 8).Value = 0.33
$ws.Cells.Item(5, 9).Value = 0.45
$ws.Cells.Item(5, 10).Value = 0.11
$ws.Cells.Item(5, 11).Value = 0.32

# Row 6
$ws.Cells.Item(6, 1).Value = 38
$ws.Cells.Item(6, 2).Value = 39
$ws.Cells.Item(6, 3).Value = 377.9573794681603
$ws.Cells.Item(6, 4).Value = 0.0555556
$ws.Cells.Item(6, 5).Value = 0.0555556
$ws.Cells.Item(6, 6).Value = 0.1111112
$ws.Cells.Item(6, 7).Value = 566.9360692022404
$ws.Cells.Item(6, 8).Value = 0.37
$ws.Cells.Item(6, 9).Value = 0.48
$ws.Cells.Item(6, 10).Value = 0.14
$ws.Cells.Item(6, 11).Value = 0.33

# Row 7
$ws.Cells.Item(7, 1).Value = 39
$ws.Cells.Item(7, 2).Value = 40
$ws.Cells.Item(7, 3).Value = 354.1473291618483
$ws.Cells.Item(7, 4).Value = 0.01111112
$ws.Cells.Item(7, 5).Value = 0.01111112
$ws.Cells.Item(7, 6).Value = 0.02222224
$ws.Cells.Item(7, 7).Value = 531.2209937427724
$ws.Cells.Item(7, 8).Value = 0.31
$ws.Cells.Item(7, 9).Value = 0.44
$ws.Cells.Item(7, 10).Value = 0.18
$ws.Cells.Item(7, 11).Value = 0.32

# Row 8
$ws.Cells.Item(8, 1).Value = 40
$ws.Cells.Item(8, 2).Value = 41
$ws.Cells.Item(8, 3).Value = 356.5382138814032
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 534.8073208221048
$ws.Cells.Item(8, 8).Value = 0.33
$ws.Cells.Item(8, 9).Value = 0.5
$ws.Cells.Item(8, 10).Value = 0.12
$ws.Cells.Item(8, 11).Value = 0.36

# Row 9
$ws.Cells.Item(9, 1).Value = 41
$ws.Cells.Item(9, 2).Value = 42
$ws.Cells.Item(9, 3).Value = 444.5860015701382
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 666.8790023552073
$ws.Cells.Item(9, 8).Value = 0.3
$ws.Cells.Item(9, 9).Value = 0.45
$ws.Cells.Item(9, 10).Value = 0.15
$ws.Cells.Item(9, 11).Value = 0.35

# Row 10
$ws.Cells.Item(10, 1).Value = 42
$ws.Cells.Item(10, 2).Value = 43
$ws.Cells.Item(10, 3).Value = 515.7790399549028
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 773.6685599323542
$ws.Cells.Item(10, 8).Value = 0.4
$ws.Cells.Item(10, 9).Value = 0.48
$ws.Cells.Item(10, 10).Value = 0.19
$ws.Cells.Item(10, 11).Value = 0.4

# Row 11
$ws.Cells.Item(11, 1).Value = 43
$ws.Cells.Item(11, 2).Value = 44
$ws.Cells.Item(11, 3).Value = 517.557383961184
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 776.336075941776
$ws.Cells.Item(11, 8).Value = 0.35
$ws.Cells.Item(11, 9).Value = 0.43
$ws.Cells.Item(11, 10).Value = 0.19
$ws.Cells.Item(11, 11).Value = 0.39

# Row 12
$ws.Cells.Item(12, 1).Value = 44
$ws.Cells.Item(12, 2).Value = 45
$ws.Cells.Item(12, 3).Value = 423.5225047846375
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 635.2837571769562
$ws.Cells.Item(12, 8).Value = 0.38
$ws.Cells.Item(12, 9).Value = 0.46
$ws.Cells.Item(12, 10).Value = 0.11
$ws.Cells.Item(12, 11).Value = 0.31

# Row 13
$ws.Cells.Item(13, 1).Value = 45
$ws.Cells.Item(13, 2).Value = 46
$ws.Cells.Item(13, 3).Value = 330.2780007219939
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 495.4170010829909
$ws.Cells.Item(13, 8).Value = 0.38
$ws.Cells.Item(13, 9).Value = 0.44
$ws.Cells.Item(13, 10).Value = 0.11
$ws.Cells.Item(13, 11).Value = 0.35

# Row 14
$ws.Cells.Item(14, 1).Value = 46
$ws.Cells.Item(14, 2).Value = 47
$ws.Cells.Item(14, 3).Value = 301.7454591101147
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 452.6181886651721
$ws.Cells.Item(14, 8).Value = 0.34
$ws.Cells.Item(14, 9).Value = 0.4
$ws.Cells.Item(14, 10).Value = 0.2
$ws.Cells.Item(14, 11).Value = 0.3

# Row 15
$ws.Cells.Item(15, 1).Value = 47
$ws.Cells.Item(15, 2).Value = 48
$ws.Cells.Item(15, 3).Value = 255.9827400151617
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 383.9741100227426
$ws.Cells.Item(15, 8).Value = 0.32
$ws.Cells.Item(15, 9).Value = 0.42
$ws.Cells.Item(15, 10).Value = 0.2
$ws.Cells.Item(15, 11).Value = 0.3

# Row 16
$ws.Cells.Item(16, 1).Value = 48
$ws.Cells.Item(16, 2).Value = 49
$ws.Cells.Item(16, 3).Value = 192.8515277922028
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 289.2772916883042
$ws.Cells.Item(16, 8).Value = 0.34
$ws.Cells.Item(16, 9).Value = 0.48
$ws.Cells.Item(16, 10).Value = 0.18
$ws.Cells.Item(16, 11).Value = 0.35

# Row 17
$ws.Cells.Item(17, 1).Value = 49
$ws.Cells.Item(17, 2).Value = 50
$ws.Cells.Item(17, 3).Value = 138.809629379121
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 208.2144440686815
$ws.Cells.Item(17, 8).Value = 0.4
$ws.Cells.Item(17, 9).Value = 0.48
$ws.Cells.Item(17, 10).Value = 0.2
$ws.Cells.Item(17, 11).Value = 0.34

